$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 81: update hours and revise the task note text
$ws.Range("B81").Value = 6
$ws.Range("D81").Value = "Reviewed Kelly's project. Mostly finished front end of Indie project."

# Row 82: fill in date, hours, and task note (previously blank except for style)
$ws.Range("A82").Value = 43593
$ws.Range("B82").Value = 1
$ws.Range("D82").Value = "Looked into why new stories weren't searchable on aws."

# Row 85: add a new task note (set before D83's text so shared-string
# insertion order matches the target: "3:45 - x" ends up before the
# "Revised search jsp..." string in the shared strings table)
$ws.Range("D85").Value = "3:45 - x"

# Row 83: fill in date and replace the task note text
$ws.Range("A83").Value = 43594
$ws.Range("D83").Value = "Revised search jsp to hold inputs when search validation failed."
